$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Programa"
